$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 previously held "Dronekit studieren" - update text to "Dronekit  studieren"
$ws.Range("B13").Value = "Dronekit  studieren"

# Row 14 previously held "GPS der Drohne auslesen" - keep same text (position swap in shared strings
# happens naturally as the underlying table is rebuilt), re-assign to be explicit.
$ws.Range("B14").Value = "GPS der Drohne auslesen"

# Fill in the two new tasks (rows 15-16), entering column by column to mirror
# how the shared string table ends up ordered in the saved file.
$ws.Range("B15").Value = "neuen Raspberry aufsetzen"
$ws.Range("B16").Value = "neuen Raspberry in Drohne einbauen"

$ws.Range("C15").Value = "Emanuel"
$ws.Range("C16").Value = "Martin / Emanuel / … (?)"

$ws.Range("D15").Value = "x"
$ws.Range("D16").Value = "x"

$ws.Range("F15").Value = "offen"
$ws.Range("F16").Value = "offen"

# Update the active selection to F17, matching the saved cursor position in the file
$ws.Range("F17").Select()
